$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 53 values that changed
$ws.Range("E53").Value = 73393
$ws.Range("G53").Value = 44060
$ws.Range("T53").Value = 177316
$ws.Range("V53").Value = 122708
$ws.Range("W53").Value = 240426

# Add new row 54 data
# A54 needs to be plain text "01-04-2021" (shared string), not an auto-converted
# date serial. Writing it as a formula that evaluates to the string, then
# converting that formula to its static value via copy/paste-special values,
# avoids Excel's automatic date recognition (and avoids creating extra
# number-format/style entries).
$ws.Range("A54").Formula = "=""01-04-2021"""
$ws.Range("A54").Copy()
$ws.Range("A54").PasteSpecial(-4163)
$ws.Range("B54").Value = 31868
$ws.Range("C54").Value = 31300
$ws.Range("D54").Value = 568
$ws.Range("E54").Value = 70629
$ws.Range("F54").Value = 26872
$ws.Range("G54").Value = 43757
$ws.Range("H54").Value = 7454
$ws.Range("I54").Value = 309
$ws.Range("J54").Value = 7145
$ws.Range("K54").Value = 24049
$ws.Range("L54").Value = 3
$ws.Range("M54").Value = 24045
$ws.Range("N54").Value = 51620
$ws.Range("O54").Value = 3711
$ws.Range("P54").Value = 47909
$ws.Range("Q54").Value = 149
$ws.Range("R54").Value = 0
$ws.Range("S54").Value = 149
$ws.Range("T54").Value = 185769
$ws.Range("U54").Value = 62195
$ws.Range("V54").Value = 123573
$ws.Range("W54").Value = 252513
